$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IndividualBiometrics")
$ws.Range("H1").Value = "Protein"
$ws.Range("I1").Value = "Ontogeny"
$ws.Range("H1:I1").Select()
